# Pulled PC13 off of JP24 and replaced it with PD7.
# This frees up PC13 for its one and only alternate function: RTC_AF1.
#
# On the "Arduino Connections" sheet, row 25 (Arduino pin 23 -> PC13) is
# removed entirely (shifting everything below it up by one row), and three
# new rows are inserted further down (just above the PC12 row) to add the
# PD7 / PD6 (U2_RX) / PD5 (U2_TX) pins.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arduino Connections")

# Remove the PC13 row (Arduino pin index 23).
$ws.Rows("25").Delete()

# After the delete, the old "PC12" row (pin 31) has shifted up to row 34.
# Insert three fresh rows right above it for the new PD7/PD6/PD5 signals.
$ws.Rows("34:36").Insert()

# Column A is just the (manually maintained) Arduino pin index, i.e. row-2.
# Re-stamp it for every row from the deletion point through the new bottom
# of the table so it stays a contiguous 0..38 sequence.
for ($r = 25; $r -le 40; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Fill in the new pin-mapping rows. Column D first (matches data-entry
# order), then the E-column alternate-function notes.
$ws.Range("D34").Value = "PD7"
$ws.Range("D35").Value = "PD6"
$ws.Range("D36").Value = "PD5"
$ws.Range("E36").Value = "U2_TX"
$ws.Range("E35").Value = "U2_RX"

# Update the remembered selections on the other two sheets that were
# visited during this editing session...
$wsDual = $wb.Worksheets.Item("XMOS Dualchip")
$wsDual.Activate()
$wsDual.Range("P15").Select()

$wsPlan = $wb.Worksheets.Item("XMOS Dualchip Planning")
$wsPlan.Activate()
$wsPlan.Range("I12").Select()

# ...then bring "Arduino Connections" to the front as the final active
# sheet/selection, matching where the author was last working.
$ws.Activate()
$ws.Range("F35").Select()
